# AFFINE CAPHER.xlsx - "Add files via upload"
#
# Updates the letters typed into the Ciphertext (row 16), Enkripsi (row 20)
# and Dekripsi (row 24) tables on Sheet1, columns D..R, and refreshes the
# active cell selection / scroll position used when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 "Ciphertext" and Row 20 "Enkripsi " share the same new letters.
$newCipherLetters = @{
    "D" = "R"; "E" = "X"; "F" = "X"; "G" = "F"; "H" = "E"
    "I" = "F"; "J" = "M"; "K" = "D"; "L" = "Z"; "M" = "F"
    "N" = "O"; "O" = "T"; "P" = "S"; "Q" = "X"; "R" = "F"
}
foreach ($col in $newCipherLetters.Keys) {
    $ws.Range($col + "16").Value = $newCipherLetters[$col]
    $ws.Range($col + "20").Value = $newCipherLetters[$col]
}

# Row 24 "Dekripsi" is updated so the decrypted text matches the original
# plaintext in row 12.
$newPlainLetters = @{
    "D" = "E"; "E" = "G"; "F" = "G"; "G" = "A"; "H" = "R"
    "I" = "A"; "J" = "L"; "K" = "I"; "L" = "Y"; "M" = "A"
    "N" = "D"; "O" = "W"; "P" = "N"; "Q" = "G"; "R" = "A"
}
foreach ($col in $newPlainLetters.Keys) {
    $ws.Range($col + "24").Value = $newPlainLetters[$col]
}

# Restore the selection / scroll state that was active when the workbook
# was saved.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H9").Select() | Out-Null
